# Apply the cryptos-list price/volume update described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.575.45"
$ws.Range("E2").Value = "  +3.54%  "

$ws.Range("D3").Value = "2.253.98"
$ws.Range("E3").Value = "  +1.89%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.528"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.08%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  +1.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "31.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.30%  "

$ws.Range("E11").Value = "  +2.81%  "

$ws.Range("E12").Value = "  +1.64%  "

$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("E14").Value = "  +2.57%  "

$ws.Range("D15").Value = "2.606.15"
$ws.Range("E15").Value = "  +2.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("D17").Value = "2.255.90"
$ws.Range("E17").Value = "  +3.99%  "

$ws.Range("E18").Value = "  +2.81%  "

$ws.Range("D19").Value = "41.502.30"
$ws.Range("E19").Value = "  +3.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.69%  "

$ws.Range("E21").Value = "  +1.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.85%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.97%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("E27").Value = "  +4.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.69%  "

$ws.Range("E30").Value = "  -0.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0739"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.31%  "

$ws.Range("E36").Value = "  -1.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.99%  "

$ws.Range("E38").Value = "  +2.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.10%  "

$ws.Range("E40").Value = "  +3.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.63%  "

$ws.Range("D43").Value = "2.050.16"
$ws.Range("E43").Value = "  -1.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.73%  "

$ws.Range("E45").Value = "  +1.94%  "

$ws.Range("E46").Value = "  +2.14%  "

$ws.Range("E47").Value = "  +5.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.99%  "

$ws.Range("E49").Value = "  +3.66%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.26%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.39%  "

